# Generate Report for Handback
#
# This localization-status workbook tracks, per target locale sheet
# (zh-cn, de-de), the handoff/handback lifecycle of each source file.
# A handback just completed for both locales, so:
#   - the status text moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is shown
#     (Overview!E:F and the per-locale Status column)
#   - the per-locale "Latest Target File" / "Latest Handback File"
#     columns (I/J) get populated with the handed-back file names,
#     with I also becoming a hyperlink (like column A already is)
#   - the per-locale "Latest Handback DateTime" column (K) is stamped
#     with the handback time (zh-cn and de-de finished at different
#     times)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$wsOverview.Range("E2").Value2 = $statusText
$wsOverview.Range("F2").Value2 = $statusText
$wsOverview.Range("E3").Value2 = $statusText
$wsOverview.Range("F3").Value2 = $statusText

# --- Per file-row hyperlink targets (same markdown files as column A) ---
$acbUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60b9f927591dcf76bb6b6345f2dc81006378ac7d/e2e/acb3d08b-601e-4505-b3a3-5b94ba208151.md"
$aeUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60b9f927591dcf76bb6b6345f2dc81006378ac7d/e2e/ae20c40a-579e-4708-88c6-9d041cfce420.md"
$acbName = "acb3d08b-601e-4505-b3a3-5b94ba208151.md"
$aeName  = "ae20c40a-579e-4708-88c6-9d041cfce420.md"

function Update-LocaleSheet($ws, $zhOrDeSuffix, $handback2, $handback3) {
    # Status column (C) on this locale sheet mirrors the Overview text.
    $ws.Range("C2").Value2 = $statusText
    $ws.Range("C3").Value2 = $statusText

    # Latest Target File (I) / Latest Handback File (J) for both rows.
    $ws.Range("I2").Value2 = $acbName
    $ws.Range("J2").Value2 = "acb3d08b-601e-4505-b3a3-5b94ba208151.dafc2fc1a905c76ab6110551a560695cdf4f3527.$zhOrDeSuffix.xlf"

    $ws.Range("I3").Value2 = $aeName
    $ws.Range("J3").Value2 = "ae20c40a-579e-4708-88c6-9d041cfce420.c9cd2fe18d1aaefff2835fed38c90737da6cd0fd.$zhOrDeSuffix.xlf"

    # Latest Handback DateTime (K) for both rows.
    $ws.Range("K2").Value2 = $handback2
    $ws.Range("K3").Value2 = $handback3

    # Match the look of the existing column-A hyperlinks (column I now
    # carries the same kind of link). Rebuild the hyperlink collection in
    # left-to-right / top-to-bottom order so relationship ids line up the
    # way Excel would emit them.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $acbUrl, "", "", $acbName)
    $ws.Hyperlinks.Add($ws.Range("I2"), $acbUrl, "", "", $acbName)
    $ws.Hyperlinks.Add($ws.Range("A3"), $aeUrl, "", "", $aeName)
    $ws.Hyperlinks.Add($ws.Range("I3"), $aeUrl, "", "", $aeName)

    # Give the new I-column links the same visual hyperlink styling as
    # column A (underlined, hyperlink-blue).
    $ws.Range("I2").Font.Underline = $true
    $ws.Range("I2").Font.Color = 15570276
    $ws.Range("I3").Font.Underline = $true
    $ws.Range("I3").Font.Color = 15570276

    # The Status/Target/Handback columns now hold much longer strings;
    # widen them the same way Excel's column autosize would.
    $ws.Columns.Item(3).ColumnWidth = 29.14
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}

Update-LocaleSheet $wsZhCn "zh-cn" "2016-09-05 05:08:06" "2016-09-05 05:08:06"
Update-LocaleSheet $wsDeDe "de-de" "2016-09-05 05:08:17" "2016-09-05 05:08:17"

# Overview's zh-cn/de-de columns (E, F) also widen now that the status
# text is longer.
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14
